# This script updates the "York (NE)_B" team-specific transition-probability
# matrix on Sheet1 of the active workbook. Rows/cols B:S (2..19) represent a
# Markov-style transition matrix whose rows must sum to 1. Several rows that
# were previously all-zero (meaning "not yet simulated") now contain the
# simulated transition probabilities following the addition of more games,
# a faster simulate-game routine, and the first pass of the optimization
# logic referenced in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.3
$ws.Range("C2").Value = 0.4
$ws.Range("P2").Value = 0.2
$ws.Range("S2").Value = 0.1

# Row 3
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.25

# Row 4
$ws.Range("P4").Value = 1

# Row 6
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.1666666666666667
$ws.Range("S6").Value = 0.6666666666666666

# Row 7
$ws.Range("Q7").Value = 0.2
$ws.Range("S7").Value = 0.8

# Row 8
$ws.Range("B8").Value = 0.12
$ws.Range("F8").Value = 0.04
$ws.Range("J8").Value = 0.08
$ws.Range("Q8").Value = 0.16
$ws.Range("R8").Value = 0.16
$ws.Range("S8").Value = 0.44

# Row 9
$ws.Range("B9").Value = 0.25
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("S9").Value = 0.5833333333333334

# Row 10
$ws.Range("B10").Value = 0.05
$ws.Range("D10").Value = 0.025
$ws.Range("F10").Value = 0.075
$ws.Range("J10").Value = 0.15
$ws.Range("O10").Value = 0.025
$ws.Range("Q10").Value = 0.2
$ws.Range("R10").Value = 0.05
$ws.Range("S10").Value = 0.425

# Row 11
$ws.Range("G11").Value = 0.1428571428571428
$ws.Range("K11").Value = 0.1428571428571428
$ws.Range("L11").Value = 0.7142857142857143

# Row 12
$ws.Range("G12").Value = 0.4
$ws.Range("J12").Value = 0.2
$ws.Range("S12").Value = 0.4

# Row 13
$ws.Range("G13").Value = 1

# Row 14
$ws.Range("G14").Value = 1

# Row 15
$ws.Range("H15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.1666666666666667
$ws.Range("O15").Value = 0.1666666666666667

# Row 16
$ws.Range("I16").Value = 0.3333333333333333
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("K16").Value = 0.1666666666666667
$ws.Range("S16").Value = 0.1666666666666667

# Row 17
$ws.Range("H17").Value = 0.3125
$ws.Range("I17").Value = 0.1875
$ws.Range("J17").Value = 0.375
$ws.Range("K17").Value = 0.0625
$ws.Range("O17").Value = 0.0625

# Row 18
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.5714285714285714
$ws.Range("K18").Value = 0.1428571428571428
$ws.Range("S18").Value = 0.1428571428571428

# Row 19
$ws.Range("H19").Value = 0.2962962962962963
$ws.Range("I19").Value = 0.1296296296296296
$ws.Range("J19").Value = 0.3703703703703703
$ws.Range("K19").Value = 0.03703703703703703
$ws.Range("N19").Value = 0.01851851851851852
$ws.Range("O19").Value = 0.05555555555555555
$ws.Range("S19").Value = 0.09259259259259259

$wb.Save()
